# docs(wbs): refresh milestone and row execution statuses
#
# The "Status Updated On" column (I) on the WBS sheet is a date serial
# (formatted yyyy-mm-dd via style 51). Every populated row (2 through 137)
# currently shows 46077 (2026-02-24) and needs to roll forward one day to
# 46078 (2026-02-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WBS")

$ws.Range("I2:I137").Value = 46078
